$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.205.84'
$ws.Range('E2').Value = '  +1.91%  '
$ws.Range('D3').Value = '2.022.78'
$ws.Range('E3').Value = '  +3.48%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '246.80'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +1.52%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.627'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +0.02%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '60.48'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +0.24%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.389'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +2.98%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0809'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +2.41%  '
$ws.Range('E11').Value = '  +1.28%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '15.05'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +6.26%  '
$ws.Range('D13').Value = '2.322.59'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.851'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +1.59%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '21.89'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +1.59%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '5.45'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +3.42%  '
$ws.Range('D17').Value = '2.016.34'
$ws.Range('E17').Value = '  +2.82%  '
$ws.Range('D18').Value = '37.174.16'
$ws.Range('E18').Value = '  +1.88%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '70.38'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +1.63%  '
$ws.Range('D20').Value = '0.0₃0864'
$ws.Range('E20').Value = '  +1.24%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '5.23'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +2.96%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '230.68'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +0.55%  '
$ws.Range('E23').Value = '  -0.05%  '
$ws.Range('E24').Value = '  +4.49%  '
$ws.Range('E25').Value = '  -0.74%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '9.38'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +2.36%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '163.74'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +1.94%  '
$ws.Range('E28').Value = '  -3.26%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '19.79'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +2.53%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.39'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +6.54%  '
$ws.Range('E31').Value = '  +1.05%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.0674'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +9.93%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '4.77'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +0.00%  '
$ws.Range('E34').Value = '  +10.86%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '4.46'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +0.39%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '3.63'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +5.25%  '
$ws.Range('E37').Value = '  +0.05%  '
$ws.Range('E38').Value = '  +1.50%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '5.33'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -2.18%  '
$ws.Range('E40').Value = '  +3.34%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.0975'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +0.82%  '
$ws.Range('E42').Value = '  +2.61%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '16.91'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +6.52%  '
$ws.Range('E44').Value = '  +1.25%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '91.36'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +2.92%  '
$ws.Range('D46').Value = '1.379.45'
$ws.Range('E46').Value = '  +1.29%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '1.05'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +2.77%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '7.44'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +3.97%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.10'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +14.08%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '46.22'
$ws.Range('D51').Style = "Normal"
